# Update the workflow file paths for the refreshed NVIS Extant and Pre1750
# intermediate rasters (date stamp 20240730/20240709 -> 20240801), as
# described in the commit message "Updated all workflow files with new
# path to NVIS Extant and Pre1750 files".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# B4 holds the RawDataPath for the Terrestrial-Pre-IUCNGET dataset.
$ws.Range("B4").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_PRE1750_IUCNGET_DK_20240801.tif"

# B5 holds the RawDataPath for the Terrestrial-Extant-IUCNGET dataset.
$ws.Range("B5").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_IUCNGET_DK_20240801.tif"

# Reflect the author's final cursor position/selection on the sheet.
$ws.Activate()
$ws.Range("B6").Select()
